$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.736.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -6.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.333.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.334.31'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.73%  '
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.118'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.375'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.900.85'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.53%  '
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.331.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000167'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.66%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '59.967.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('E21').Value = '  -8.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '353.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.556'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.464.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '68.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.47%  '
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.44'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.49%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.52'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.78%  '
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('E32').Value = '  -4.38%  '
$ws.Range('E33').Value = '  -2.04%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.369.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E37').Value = '  +2.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '157.88'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0756'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('E44').Value = '  +6.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.747'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.20'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('E48').Value = '  -3.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.27%  '
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.41'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +15.82%  '
